$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet has columns: A=code, B=status, C=group-name, D=group-code.
# This edit swaps columns C and D so that C=group-code, D=group-name,
# for the header row and every data row.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value()
    $dVal = $dCell.Value()
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
